$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update status-check timestamp in F1
$ws.Range("F1").Value = "Last status check on: 10.01.2022 08:07"

# Convert D10 from text "-0.3" to a real number
$ws.Range("D10").Value = -0.3

# Convert E10 from text date to a real Excel date/time serial value,
# matching the numeric date formatting used by the other rows (column E)
$ws.Range("E10").Value = Get-Date -Year 2022 -Month 1 -Day 10 -Hour 8 -Minute 4 -Second 35
$ws.Range("E10").NumberFormat = $ws.Range("E9").NumberFormat
